# Add a new "November_2025" worksheet as the first tab, based on the
# existing "October_2025" sheet, with updated title + data values.

$wb = $excel.ActiveWorkbook

$october = $wb.Worksheets.Item("October_2025")

# Copy October_2025 to before itself -> new copy becomes the first sheet,
# pushing October_2025 (and everything else) one position later.
$october.Copy($october)

$november = $wb.Worksheets.Item(1)
$november.Name = "November_2025"

$a1Style = $november.Range("A1").Style
$november.Range("A1").Value = "Table J (11/28/2025) Swap Benchmark Spreads (in bps)"
$november.Range("A1").Style = $a1Style

$november.Range("B4").Value = -1.887906877
$november.Range("C4").Value = -8.6326097270000002

$november.Range("B5").Value = 1.7783312499999999
$november.Range("C5").Value = -19.17339213

$november.Range("B6").Value = -6.8667596890000002
$november.Range("C6").Value = -6.6292665910000004

$november.Range("B7").Value = -14.837141600000001
$november.Range("C7").Value = -9.5367021649999995

$november.Range("B8").Value = -20.829166480000001
$november.Range("C8").Value = -13.69393036

$november.Range("B9").Value = -22.76229987
$november.Range("C9").Value = -17.806993559999999

$november.Range("B10").Value = -24.877004639999999
$november.Range("C10").Value = -20.803627909999999

$november.Range("B11").Value = -28.76618663
$november.Range("C11").Value = -23.199092140000001

$november.Range("B12").Value = -32.13687204
$november.Range("C12").Value = -24.827642269999998

$november.Range("B13").Value = -34.36064159
$november.Range("C13").Value = -25.966189230000001

$november.Range("B14").Value = -36.217203499999997
$november.Range("C14").Value = -26.83723874

$november.Range("B15").Value = -37.636845530000002
$november.Range("C15").Value = -27.58993267

$november.Range("B16").Value = -39.932199699999998
$november.Range("C16").Value = -28.524672379999998

$november.Range("B17").Value = -43.440759329999999
$november.Range("C17").Value = -29.76895309

$november.Range("B18").Value = -47.18283701
$november.Range("C18").Value = -31.36868042

$november.Range("B19").Value = -50.577111940000002
$november.Range("C19").Value = -33.190270419999997

$november.Range("B20").Value = -53.492839099999998
$november.Range("C20").Value = -35.153061600000001

$november.Range("B21").Value = -56.023917300000001
$november.Range("C21").Value = -37.218337740000003

$november.Range("B22").Value = -58.055090569999997
$november.Range("C22").Value = -39.341179740000001

$november.Range("B23").Value = -59.645390429999999
$november.Range("C23").Value = -41.460715399999998

$november.Range("B24").Value = -60.875489780000002
$november.Range("C24").Value = -43.518742189999998

$november.Range("B25").Value = -61.792998410000003
$november.Range("C25").Value = -45.453901119999998

$november.Range("B26").Value = -62.507469839999999
$november.Range("C26").Value = -47.303131319999999

$november.Range("B27").Value = -63.138604659999999
$november.Range("C27").Value = -49.077808040000001

$november.Range("B28").Value = -63.692138669999999
$november.Range("C28").Value = -50.70311641

$november.Range("B29").Value = -64.249029350000001
$november.Range("C29").Value = -52.115476839999999

$november.Range("B30").Value = -64.798857859999998
$november.Range("C30").Value = -53.254990900000003

$november.Range("B31").Value = -65.321158749999995
$november.Range("C31").Value = -54.088753140000001

$november.Range("B32").Value = -65.835949240000005
$november.Range("C32").Value = -54.630743340000002

$november.Range("B33").Value = -66.429597860000001
$november.Range("C33").Value = -54.882454289999998

$november.Range("B34").Value = -67.028187279999997
$november.Range("C34").Value = -54.841728379999999

$november.Range("B35").Value = -67.757695569999996
$november.Range("C35").Value = -54.54834915

$november.Range("B36").Value = -45.036781870187504
$november.Range("C36").Value = -34.848177606343747
